# This workbook's weekly data refresh permutes the per-record data blocks
# (Fecha, Variedad, Volumen, Precio minimo/maximo/promedio, Unidad de
# comercializacion, Origen, Precio $/Kg, Kg o Unidades) across the existing
# rows 2-62, while the constant columns (Mercado ID, Mercado, Region,
# Codreg, Categoria ID, Categoria, Calidad, Clasificacion) stay untouched.
#
# We read every row's mutable fields live via COM (.Value2, which returns
# plain scalars instead of date/variant wrapper objects in this runtime),
# then write them back out in the permuted order below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 62

# Columns that move together as one record, keyed by column letter -> index.
$cols = @{
    D = 4   # Fecha
    H = 8   # Variedad
    J = 10  # Volumen
    K = 11  # Precio minimo
    L = 12  # Precio maximo
    M = 13  # Precio promedio ponderado
    N = 14  # Unidad de comercializacion
    O = 15  # Origen
    P = 16  # Precio $/Kg
    Q = 17  # Kg o Unidades
}

# Snapshot the current ("before") values of the mutable columns for every row.
$rows = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $row = @{
        D = $ws.Cells.Item($r, $cols.D).Value2
        H = $ws.Cells.Item($r, $cols.H).Value2
        J = $ws.Cells.Item($r, $cols.J).Value2
        K = $ws.Cells.Item($r, $cols.K).Value2
        L = $ws.Cells.Item($r, $cols.L).Value2
        M = $ws.Cells.Item($r, $cols.M).Value2
        N = $ws.Cells.Item($r, $cols.N).Value2
        O = $ws.Cells.Item($r, $cols.O).Value2
        P = $ws.Cells.Item($r, $cols.P).Value2
        Q = $ws.Cells.Item($r, $cols.Q).Value2
    }
    $rows += $row
}

# For target row (index 0 => row 2, ... index 60 => row 62), this gives the
# sheet row number (2-62) whose "before" record now belongs there.
$mapping = @(
    15, 2, 54, 3, 12, 57, 31, 51, 44, 28,
    52, 48, 40, 25, 26, 58, 59, 33, 10, 56,
    35, 42, 47, 55, 29, 21, 22, 24, 6, 34,
    30, 49, 43, 27, 50, 53, 13, 45, 14, 5,
    46, 23, 19, 20, 37, 9, 7, 61, 11, 62,
    17, 4, 39, 60, 32, 36, 41, 18, 8, 16,
    38
)

for ($i = 0; $i -lt $mapping.Count; $i++) {
    $r = $firstRow + $i
    $srcRow = $mapping[$i]
    $src = $rows[$srcRow - $firstRow]

    $ws.Cells.Item($r, $cols.D).Value2 = $src.D
    $ws.Cells.Item($r, $cols.H).Value2 = $src.H
    $ws.Cells.Item($r, $cols.J).Value2 = $src.J
    $ws.Cells.Item($r, $cols.K).Value2 = $src.K
    $ws.Cells.Item($r, $cols.L).Value2 = $src.L
    $ws.Cells.Item($r, $cols.M).Value2 = $src.M
    $ws.Cells.Item($r, $cols.N).Value2 = $src.N
    $ws.Cells.Item($r, $cols.O).Value2 = $src.O
    $ws.Cells.Item($r, $cols.P).Value2 = $src.P
    $ws.Cells.Item($r, $cols.Q).Value2 = $src.Q
}
